# Add the 2023 column (S) to the Adjaria AR hotels & restaurants indicator
# table: "add genders in batumi" -- append the newest reporting year's
# figures next to the existing 2006-2022 (B:R) series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New year header + the twelve indicator rows for 2023, one per data row
# (mirrors the layout of the existing B:R year columns).
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 455.5
$ws.Range("S5").Value = 483.8
$ws.Range("S6").Value = 5833
$ws.Range("S7").Value = 5749
$ws.Range("S8").Value = 1408.7
$ws.Range("S9").Value = 253.6
$ws.Range("S10").Value = 99.1
$ws.Range("S11").Value = 230.1
$ws.Range("S12").Value = 233.6
$ws.Range("S13").Value = 260.3
$ws.Range("S14").Value = 0.1

# Mirror the look of the preceding (2022 / column R) cells onto the new
# 2023 / column S cells -- same number format, font, borders and alignment
# for every row in the table.
$ws.Range("R3:R14").Copy() | Out-Null
$ws.Range("S3:S14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the selection where the author left it after the edit.
$ws.Range("B8").Select() | Out-Null
